# Översikt FÄRGELANDA.xlsx update
# - Column C ("Förändrad") date bumped from 2023-09-03 (45172) to 2023-09-06 (45175) for every data row (2..319)
# - Row 3 (A 33191-2022) gained a new signal species "Gul taggsvamp":
#     NT (J3) 1 -> 2, Rödlistade (O3) 2 -> 3, Alla arter (Q3) 4 -> 5
#     Artnamn (R3) gets "Gul taggsvamp" inserted as the second line

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column for all data rows at once.
$ws.Range("C2:C319").Value = 45175

# Row-specific updates for A 33191-2022 (row 3)
$ws.Range("J3").Value = 2
$ws.Range("O3").Value = 3
$ws.Range("Q3").Value = 5
$ws.Range("R3").Value = "Knärot`r`nGul taggsvamp`r`nSkirmossa`r`nHavstulpanlav`r`nKorallblylav"

# Writing to the wrap-text R3 cell triggers this engine's auto row-height
# recompute; restore the original (unchanged) row height afterwards.
$ws.Rows.Item(3).RowHeight = 15
